# Article details page request finished
# Update the "tags" column (G) on the Articles sheet with the new values
# returned by the finished request.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Articles")

$ws.Range("G2").Value = "1.1.list, 2.2of, 3.3.nsmi, tags"
$ws.Range("G3").Value = "1.3.of, 2.3.nsmi, tags"
$ws.Range("G4").Value = "list, tags"

$ws.Activate()
$ws.Range("I3").Select()
